# Insert a new weekly price row at row 220 (pushing existing rows 220-228
# down to 221-229), then populate the new row with the latest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 220:228 down to 221:229 to make room for the new record.
$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new weekly record.
$ws.Cells.Item(220, 1).Value = 3
$ws.Cells.Item(220, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(220, 3).Value = "Coquimbo"
$ws.Cells.Item(220, 4).Value = 44509
$ws.Cells.Item(220, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(220, 5).Value = 5
$ws.Cells.Item(220, 6).Value = 100112043
$ws.Cells.Item(220, 7).Value = "Pepino ensalada"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 130
$ws.Cells.Item(220, 11).Value = 7000
$ws.Cells.Item(220, 12).Value = 7500
$ws.Cells.Item(220, 13).Value = 7269
$ws.Cells.Item(220, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(220, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(220, 16).Value = 104
$ws.Cells.Item(220, 17).Value = 70
$ws.Cells.Item(220, 18).Value = "Hortaliza"
